$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.856.10'
$ws.Range('E2').Value = '  -4.42%  '
$ws.Range('D3').Value = '2.324.79'
$ws.Range('E3').Value = '  -5.75%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.55'
$ws.Range('E5').Value = '  -4.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '84.76'
$ws.Range('E6').Value = '  -7.49%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.530'
$ws.Range('E7').Value = '  -3.45%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.482'
$ws.Range('E9').Value = '  -4.67%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0815'
$ws.Range('E10').Value = '  -3.97%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '29.64'
$ws.Range('E11').Value = '  -9.26%  '
$ws.Range('E12').Value = '  +0.41%  '
$ws.Range('D13').Value = '2.674.77'
$ws.Range('E13').Value = '  -5.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.37'
$ws.Range('E14').Value = '  -6.79%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.79'
$ws.Range('E15').Value = '  -4.69%  '
$ws.Range('D16').Value = '2.333.72'
$ws.Range('E16').Value = '  -5.53%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.751'
$ws.Range('E17').Value = '  -4.97%  '
$ws.Range('D18').Value = '39.807.12'
$ws.Range('E18').Value = '  -4.17%  '
$ws.Range('D19').Value = '0.0₃0896'
$ws.Range('E19').Value = '  -4.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.08'
$ws.Range('E20').Value = '  -4.73%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '67.83'
$ws.Range('E21').Value = '  -6.29%  '
$ws.Range('E22').Value = '  -5.98%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.28'
$ws.Range('E23').Value = '  -0.25%  '
$ws.Range('E24').Value = '  -7.91%  '
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.78'
$ws.Range('E26').Value = '  -7.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.12'
$ws.Range('E27').Value = '  -5.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.18'
$ws.Range('E28').Value = '  -2.53%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.16'
$ws.Range('E29').Value = '  -5.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '33.79'
$ws.Range('E30').Value = '  -6.85%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '151.51'
$ws.Range('E31').Value = '  -4.84%  '
$ws.Range('E32').Value = '  +0.00%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.07'
$ws.Range('E33').Value = '  -5.84%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.44'
$ws.Range('E34').Value = '  -4.88%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0711'
$ws.Range('E35').Value = '  -5.79%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0983'
$ws.Range('E37').Value = '  -4.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.71'
$ws.Range('E38').Value = '  -6.57%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '15.33'
$ws.Range('E39').Value = '  -9.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.69'
$ws.Range('E40').Value = '  -7.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.78'
$ws.Range('E41').Value = '  -4.71%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '1.936.81'
$ws.Range('E42').Value = '  -2.59%  '
$ws.Range('B43').Value = 'ApeXProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.22'
$ws.Range('E43').Value = '  -4.57%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0260'
$ws.Range('E44').Value = '  -7.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '17.23'
$ws.Range('E45').Value = '  -7.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.34'
$ws.Range('E46').Value = '  -3.59%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.67'
$ws.Range('E47').Value = '  -9.22%  '
$ws.Range('D48').Value = '2.570.19'
$ws.Range('E48').Value = '  -4.70%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '91.87'
$ws.Range('E49').Value = '  -5.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '70.61'
$ws.Range('E50').Value = '  -5.52%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '49.67'
$ws.Range('E51').Value = '  -3.75%  '
